$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 100.833336
$ws.Range("I9").Value = 116.25
$ws.Range("K9").Value = 116.25
$ws.Range("M9").Value = 52.75

$ws.Range("H38").Value = 51.333332
$ws.Range("I38").Value = 51.333332
$ws.Range("K38").Value = 153.999996
$ws.Range("M38").Value = 218.000004

$ws.Range("H39").Value = 4466
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4466
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = 13398
$ws.Range("N39").Value = -13990

$ws.Range("H76").Value = 5612.091
$ws.Range("I76").Value = 5005.3335
$ws.Range("K76").Value = 5005.3335
$ws.Range("M76").Value = -4690.3335

$ws.Range("H79").Value = 5612.091
$ws.Range("I79").Value = 5005.3335
$ws.Range("K79").Value = 5005.3335
$ws.Range("M79").Value = -3913.3335

$ws.Range("H125").Value = 1148.3077
$ws.Range("I125").Value = 559.6667
$ws.Range("K125").Value = 5037.0003
$ws.Range("M125").Value = -2577.0003

$ws.Range("H132").Value = 3368.255
$ws.Range("I132").Value = 2525.0454
$ws.Range("J132").Value = 8668.429
$ws.Range("K132").Value = 7575.1362
$ws.Range("L132").Value = 26005.287
$ws.Range("M132").Value = -5045.1362
$ws.Range("N132").Value = -31065.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8501.35
$ws.Range("I32").Value = 5001.66
$ws.Range("K32").Value = 5001.66
$ws.Range("M32").Value = -4714.66

$ws.Range("H63").Value = 1764.0588
$ws.Range("I63").Value = 1764.0588
$ws.Range("K63").Value = 1764.0588
$ws.Range("M63").Value = -1078.0588

$ws.Range("H66").Value = 1764.0588
$ws.Range("I66").Value = 1764.0588
$ws.Range("K66").Value = 8820.294
$ws.Range("M66").Value = -5388.294

$ws.Range("H122").Value = 3595.3333
$ws.Range("I122").Value = 3419.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 10259.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7809.25
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 33717.168
$ws.Range("J76").Value = 33717.168
$ws.Range("L76").Value = 33717.168
$ws.Range("N76").Value = -34347.168

$ws.Range("H79").Value = 33717.168
$ws.Range("J79").Value = 33717.168
$ws.Range("L79").Value = 33717.168
$ws.Range("N79").Value = -35901.168

$ws.Range("H99").Value = 1958.7273
$ws.Range("I99").Value = 1569.6111
$ws.Range("J99").Value = 3709.75
$ws.Range("K99").Value = 1569.6111
$ws.Range("L99").Value = 3709.75
$ws.Range("M99").Value = -71.61110000000008
$ws.Range("N99").Value = -6705.75

$ws.Range("H105").Value = 2755.7144
$ws.Range("I105").Value = 2308.2
$ws.Range("J105").Value = 3874.5
$ws.Range("K105").Value = 2308.2
$ws.Range("L105").Value = 3874.5
$ws.Range("M105").Value = -561.1999999999998
$ws.Range("N105").Value = -7368.5

$ws.Range("H107").Value = 1014.125
$ws.Range("I107").Value = 844.7143
$ws.Range("K107").Value = 844.7143
$ws.Range("M107").Value = 1075.2857

$ws.Range("H134").Value = 1935.4468
$ws.Range("I134").Value = 1971.0217
$ws.Range("J134").Value = 299
$ws.Range("K134").Value = 5913.0651
$ws.Range("L134").Value = 897
$ws.Range("M134").Value = -3378.0651
$ws.Range("N134").Value = -5967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2591.2856
$ws.Range("I31").Value = 1211.1936
$ws.Range("J31").Value = 6480.636
$ws.Range("K31").Value = 1211.1936
$ws.Range("L31").Value = 6480.636
$ws.Range("M31").Value = -916.1936000000001
$ws.Range("N31").Value = -7070.636

$ws.Range("H34").Value = 2591.2856
$ws.Range("I34").Value = 1211.1936
$ws.Range("J34").Value = 6480.636
$ws.Range("K34").Value = 1211.1936
$ws.Range("L34").Value = 6480.636
$ws.Range("M34").Value = -1009.1936
$ws.Range("N34").Value = -6884.636

$ws.Range("H99").Value = 3785.7778
$ws.Range("I99").Value = 2292
$ws.Range("J99").Value = 4532.6665
$ws.Range("K99").Value = 2292
$ws.Range("L99").Value = 4532.6665
$ws.Range("M99").Value = -794
$ws.Range("N99").Value = -7528.6665

$ws.Range("H126").Value = 3785.7778
$ws.Range("I126").Value = 2292
$ws.Range("J126").Value = 4532.6665
$ws.Range("K126").Value = 6876
$ws.Range("L126").Value = 13597.9995
$ws.Range("M126").Value = -4406
$ws.Range("N126").Value = -18537.9995

$ws.Range("H134").Value = 44803.168
$ws.Range("I134").Value = 48628.5
$ws.Range("K134").Value = 145885.5
$ws.Range("M134").Value = -143350.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1262
$ws.Range("I5").Value = 1392.8572
$ws.Range("J5").Value = 1109.3334
$ws.Range("K5").Value = 4178.571599999999
$ws.Range("L5").Value = 3328.0002
$ws.Range("M5").Value = -4066.571599999999
$ws.Range("N5").Value = -3552.0002

$ws.Range("H37").Value = 67125.28999999999
$ws.Range("J37").Value = 67125.28999999999
$ws.Range("L37").Value = 201375.87
$ws.Range("N37").Value = -201599.87

$ws.Range("H39").Value = 4783.5
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 5226.1113
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 15678.3339
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -16266.3339

$ws.Range("H46").Value = 1337.75
$ws.Range("I46").Value = 2125
$ws.Range("J46").Value = 550.5
$ws.Range("K46").Value = 6375
$ws.Range("L46").Value = 1651.5
$ws.Range("M46").Value = -6284
$ws.Range("N46").Value = -1833.5

$ws.Range("H75").Value = 1008
$ws.Range("J75").Value = 1094.8572
$ws.Range("L75").Value = 3284.5716
$ws.Range("N75").Value = -5280.571599999999

$ws.Range("H78").Value = 1008
$ws.Range("J78").Value = 1094.8572
$ws.Range("L78").Value = 9853.7148
$ws.Range("N78").Value = -19837.7148

$ws.Range("H122").Value = 484.20834
$ws.Range("J122").Value = 583.4545000000001
$ws.Range("L122").Value = 5251.0905
$ws.Range("N122").Value = -10151.0905

$ws.Range("H135").Value = 1262
$ws.Range("I135").Value = 1392.8572
$ws.Range("J135").Value = 1109.3334
$ws.Range("K135").Value = 12535.7148
$ws.Range("L135").Value = 9984.000599999999
$ws.Range("M135").Value = -10000.7148
$ws.Range("N135").Value = -15054.0006

$ws.Range("H140").Value = 2896.15
$ws.Range("I140").Value = 2896.15
$ws.Range("K140").Value = 8688.450000000001
$ws.Range("M140").Value = -3508.450000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4998.25
$ws.Range("I80").Value = 4999
$ws.Range("J80").Value = 4998
$ws.Range("K80").Value = 4999
$ws.Range("L80").Value = 4998
$ws.Range("M80").Value = -4001
$ws.Range("N80").Value = -6994

$ws.Range("H83").Value = 4998.25
$ws.Range("I83").Value = 4999
$ws.Range("J83").Value = 4998
$ws.Range("K83").Value = 24995
$ws.Range("L83").Value = 24990
$ws.Range("M83").Value = -20003
$ws.Range("N83").Value = -34974

$ws.Range("H97").Value = 782.45
$ws.Range("I97").Value = 718.6842
$ws.Range("J97").Value = 1994
$ws.Range("K97").Value = 718.6842
$ws.Range("L97").Value = 1994
$ws.Range("M97").Value = -222.6842
$ws.Range("N97").Value = -2986

$ws.Range("H102").Value = 4728
$ws.Range("I102").Value = 2970.6667
$ws.Range("K102").Value = 2970.6667
$ws.Range("M102").Value = -1348.6667

$ws.Range("H122").Value = 2463.9
$ws.Range("I122").Value = 1968.625
$ws.Range("J122").Value = 4445
$ws.Range("K122").Value = 5905.875
$ws.Range("L122").Value = 13335
$ws.Range("M122").Value = -3455.875
$ws.Range("N122").Value = -18235

$ws.Range("H126").Value = 4814.1577
$ws.Range("I126").Value = 4089.0732
$ws.Range("K126").Value = 12267.2196
$ws.Range("M126").Value = -9797.2196

$ws.Range("H132").Value = 32195.117
$ws.Range("I132").Value = 45523.13
$ws.Range("K132").Value = 136569.39
$ws.Range("M132").Value = -134039.39

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9014.519
$ws.Range("I7").Value = 10200
$ws.Range("J7").Value = 4865.3335
$ws.Range("K7").Value = 10200
$ws.Range("L7").Value = 4865.3335
$ws.Range("M7").Value = -10088
$ws.Range("N7").Value = -5089.3335

$ws.Range("H46").Value = 15396
$ws.Range("I46").Value = 19662.766
$ws.Range("K46").Value = 19662.766
$ws.Range("M46").Value = -19474.766

$ws.Range("H126").Value = 9014.519
$ws.Range("I126").Value = 10200
$ws.Range("J126").Value = 4865.3335
$ws.Range("K126").Value = 30600
$ws.Range("L126").Value = 14596.0005
$ws.Range("M126").Value = -28130
$ws.Range("N126").Value = -19536.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 637.5833
$ws.Range("I107").Value = 340.77777
$ws.Range("K107").Value = 1022.33331
$ws.Range("M107").Value = 897.66669

$ws.Range("H132").Value = 51113.57
$ws.Range("I132").Value = 59119.723
$ws.Range("K132").Value = 177359.169
$ws.Range("M132").Value = -174829.169

$ws.Range("H136").Value = 4307.4194
$ws.Range("I136").Value = 4216.08
$ws.Range("J136").Value = 4688
$ws.Range("K136").Value = 12648.24
$ws.Range("L136").Value = 14064
$ws.Range("M136").Value = -10098.24
$ws.Range("N136").Value = -19164
